# Applies the cryptos.xlsx price/volume/ranking updates described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.912.24"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.503.26"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.88"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.32"
$ws.Range("E6").Value = "  +4.54%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("E9").Value = "  +4.16%  "
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.342"
$ws.Range("E11").Value = "  +3.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.96"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.87"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.913.71"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.721.53"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.503.71"
$ws.Range("E17").Value = "  +3.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.11"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.56"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.88"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.81"
$ws.Range("E23").Value = "  +3.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.30"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.78"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.21"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0915"
$ws.Range("E29").Value = "  +1.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "511.23"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.86"
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.27"
$ws.Range("E32").Value = "  +3.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.79"
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.124"
$ws.Range("E35").Value = "  +7.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.38"
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.46"
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.68"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.35"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.75"
$ws.Range("E41").Value = "  +3.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.331"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.89"
$ws.Range("E43").Value = "  +1.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.43"
$ws.Range("E44").Value = "  +3.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "145.91"
$ws.Range("E45").Value = "  +2.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.52"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.518"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.60"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.588"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0941"
$ws.Range("E51").Value = "  +0.51%  "
